# The author trimmed the "temperatura 0.4" observation text in cell E7
# back down to the shorter, already-used note. Shared-string table
# reindexing (dropping the now-orphaned string) is handled automatically
# by the save path, so the only real content edit is this cell's value.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resultados")

$ws.Range("E7").Value = "A resposta está correta em todos os casos."

# Cosmetic: the sheet's zoom was reset to 100% (from 90%) when the file
# was re-saved.
$excel.ActiveWindow.Zoom = 100
